# Aggiornamento al 23 agosto 2021
# Appends 14 new daily rows (344-357) to Sheet1, continuing the existing
# date series (column A) with zero values in columns B, C and D, using the
# same formatting as the last existing data row (343).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 343
$firstNewDate = 44418
$lastNewDate = 44431

# Copy the formatting (style/number format/borders) of the last data row's
# date cell (column A) down onto the new date cells; columns B-D keep the
# default (unstyled) formatting, matching the existing rows.
$ws.Range("A$lastRow").Copy()
$newFirstRow = $lastRow + 1
$newLastRow = $lastRow + ($lastNewDate - $firstNewDate + 1)
$ws.Range("A${newFirstRow}:A${newLastRow}").PasteSpecial(-4122)

$row = $newFirstRow
$dateValue = $firstNewDate
while ($row -le $newLastRow) {
    $ws.Cells.Item($row, 1).Value = $dateValue
    $ws.Cells.Item($row, 2).Value = 0
    $ws.Cells.Item($row, 3).Value = 0
    $ws.Cells.Item($row, 4).Value = 0
    $row = $row + 1
    $dateValue = $dateValue + 1
}
